$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the population / seats values for rows 19-24 (Dublin SouthWest,
# Dublin West, Dun Laoghaire, GalwayEast, GalwayWest, Kerry) using the new
# "=" style instead of the old ':"..."' quoted style.
$ws.Range("C19").Value = '{Population = 144908,'
$ws.Range("C20").Value = '{Population =113179,'
$ws.Range("D20").Value = 'Seats = 4, Name: "Dublin West"}),'
$ws.Range("C21").Value = '{Population = 118791,'
$ws.Range("C22").Value = '{Population = 89564,'
$ws.Range("C23").Value = '{Population = 150874 ,'
$ws.Range("C24").Value = '{Population = '

# Reflect the resulting selection left on the sheet after the edits -
# the whole data table (B1:D41) ends up selected.
$ws.Range("B1:D41").Select()
